$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Fecha (D) and Calidad (L) values between row 2 and row 3.
$d2 = $ws.Range("D2").Value2
$d3 = $ws.Range("D3").Value2
$ws.Range("D2").Value2 = $d3
$ws.Range("D3").Value2 = $d2

$l2 = $ws.Range("L2").Value2
$l3 = $ws.Range("L3").Value2
$ws.Range("L2").Value2 = $l3
$ws.Range("L3").Value2 = $l2
